$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title/date shared-string text (Volume number and report week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Cells whose type/style changes (text marker <-> numeric) ---
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C15").Value = 1

$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("F26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F26").Copy()
$ws.Range("C26").PasteSpecial(-4163)
$ws.Range("C26").Value = 1

$ws.Range("I30").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("I30").Copy()
$ws.Range("G30").PasteSpecial(-4163)
$ws.Range("G30").Value = 1

$ws.Range("K30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("H30").Value = -100

# --- Plain numeric value updates ---
# Row 15
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 140
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 9.090909090909

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 69
$ws.Range("J16").Value = 57
$ws.Range("K16").Value = 21.052631578947
$ws.Range("L16").Value = 40.816326530612
$ws.Range("M16").Value = 38
$ws.Range("N16").Value = -86.2

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = -50
$ws.Range("J17").Value = 81
$ws.Range("K17").Value = 3.703703703703
$ws.Range("L17").Value = 58.490566037735
$ws.Range("N17").Value = -33.333333333333

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 108
$ws.Range("K18").Value = 28.703703703703
$ws.Range("L18").Value = 6.106870229007
$ws.Range("M18").Value = 101.449275362319
$ws.Range("N18").Value = -85.165421558164

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 8.333333333333
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -4.081632653061
$ws.Range("I19").Value = 508
$ws.Range("J19").Value = 402
$ws.Range("K19").Value = 26.36815920398
$ws.Range("L19").Value = 14.932126696832
$ws.Range("M19").Value = -11.652173913043
$ws.Range("N19").Value = -70.971428571428

# Row 20
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 26.829268292682
$ws.Range("L20").Value = 40.54054054054
$ws.Range("M20").Value = 108
$ws.Range("N20").Value = -89.430894308943

# Row 21
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 10.526315789473
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -14.942528735632
$ws.Range("I21").Value = 864
$ws.Range("J21").Value = 694
$ws.Range("K21").Value = 24.495677233429
$ws.Range("L21").Value = 19.833564493758
$ws.Range("M21").Value = 12.646675358539
$ws.Range("N21").Value = -77.358490566037

# Row 22
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 25
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = 38.888888888888
$ws.Range("L22").Value = 31.578947368421
$ws.Range("M22").Value = -3.846153846153

# Row 24
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 57.142857142857
$ws.Range("I24").Value = 1014
$ws.Range("J24").Value = 758
$ws.Range("K24").Value = 33.77308707124
$ws.Range("L24").Value = 7.757704569606
$ws.Range("M24").Value = 109.504132231405

# Row 25
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 17
$ws.Range("H25").Value = -10.526315789473
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 135
$ws.Range("K25").Value = 31.851851851851
$ws.Range("L25").Value = 37.984496124031
$ws.Range("M25").Value = -2.732240437158

# Row 26
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 14
$ws.Range("K26").Value = 55.555555555555
$ws.Range("L26").Value = 7.692307692307

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 43
$ws.Range("J27").Value = 55
$ws.Range("K27").Value = -21.818181818181
$ws.Range("L27").Value = 43.333333333333

# Row 30
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 120

Write-Output "edits applied"
